$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 2
$ws.Range("H2").Value = 77.625
$ws.Range("I2").Value = 86.833336
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 86.833336
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 26.166664
$ws.Range("N2").Value = -276
# row 29
$ws.Range("H29").Value = 1657.1428
$ws.Range("I29").Value = 150
$ws.Range("J29").Value = 3666.6667
$ws.Range("K29").Value = 450
$ws.Range("L29").Value = 11000.0001
$ws.Range("M29").Value = -169
$ws.Range("N29").Value = -11562.0001
# row 38
$ws.Range("H38").Value = 106.5
$ws.Range("I38").Value = 106.5
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 319.5
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = 52.5
$ws.Range("N38").ClearContents()
# row 58
$ws.Range("H58").Value = 244.83333
$ws.Range("I58").Value = 244.83333
$ws.Range("K58").Value = 734.49999
$ws.Range("M58").Value = -584.49999
# row 70
$ws.Range("H70").Value = 1015.04346
$ws.Range("I70").Value = 824.7143
$ws.Range("J70").Value = 1311.1111
$ws.Range("K70").Value = 2474.1429
$ws.Range("L70").Value = 3933.3333
$ws.Range("M70").Value = -2204.1429
$ws.Range("N70").Value = -4473.3333
# row 73
$ws.Range("H73").Value = 1015.04346
$ws.Range("I73").Value = 824.7143
$ws.Range("J73").Value = 1311.1111
$ws.Range("K73").Value = 2474.1429
$ws.Range("L73").Value = 3933.3333
$ws.Range("M73").Value = -1538.1429
$ws.Range("N73").Value = -5805.3333
# row 87
$ws.Range("H87").Value = 39645
$ws.Range("J87").Value = 39645
$ws.Range("L87").Value = 39645
$ws.Range("N87").Value = -42141
# row 90
$ws.Range("H90").Value = 39645
$ws.Range("J90").Value = 39645
$ws.Range("L90").Value = 118935
$ws.Range("N90").Value = -131415
# row 134
$ws.Range("H134").Value = 51428.57
$ws.Range("J134").Value = 51428.57
$ws.Range("L134").Value = 51428.57
$ws.Range("N134").Value = -61568.57
# row 135
$ws.Range("H135").Value = 50016.617
$ws.Range("I135").Value = 54755.21
$ws.Range("J135").Value = 5000
$ws.Range("K135").Value = 492796.89
$ws.Range("L135").Value = 45000
$ws.Range("M135").Value = -490261.89
$ws.Range("N135").Value = -50070
# row 137
$ws.Range("H137").Value = 2441693.2
$ws.Range("I137").Value = 3335257.5
$ws.Range("J137").Value = 4700.091
$ws.Range("K137").Value = 10005772.5
$ws.Range("L137").Value = 14100.273
$ws.Range("M137").Value = -10003222.5
$ws.Range("N137").Value = -19200.273
# row 138
$ws.Range("H138").Value = 2600623.8
$ws.Range("I138").Value = 2468
$ws.Range("J138").Value = 4170342.8
$ws.Range("K138").Value = 7404
$ws.Range("L138").Value = 12511028.4
$ws.Range("M138").Value = -2264
$ws.Range("N138").Value = -12521308.4
# row 140
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
# row 141
$ws.Range("H141").Value = 7054.222
$ws.Range("I141").Value = 7498.2856
$ws.Range("J141").Value = 5500
$ws.Range("K141").Value = 22494.8568
$ws.Range("L141").Value = 16500
$ws.Range("M141").Value = -17314.8568
$ws.Range("N141").Value = -26860

$ws = $wb.Worksheets.Item("ARM")
# row 74
$ws.Range("H74").Value = 10482229
$ws.Range("I74").Value = 15922224
$ws.Range("J74").Value = 96781.82000000001
$ws.Range("K74").Value = 15922224
$ws.Range("L74").Value = 96781.82000000001
$ws.Range("M74").Value = -15921350
$ws.Range("N74").Value = -98529.82000000001
# row 77
$ws.Range("H77").Value = 10482229
$ws.Range("I77").Value = 15922224
$ws.Range("J77").Value = 96781.82000000001
$ws.Range("K77").Value = 79611120
$ws.Range("L77").Value = 483909.1
$ws.Range("M77").Value = -79606752
$ws.Range("N77").Value = -492645.1
# row 88
$ws.Range("H88").Value = 5266.75
$ws.Range("I88").Value = 2253
$ws.Range("K88").Value = 2253
$ws.Range("M88").Value = -1847
# row 91
$ws.Range("H91").Value = 5266.75
$ws.Range("I91").Value = 2253
$ws.Range("K91").Value = 2253
$ws.Range("M91").Value = -849
# row 132
$ws.Range("H132").Value = 9475332
$ws.Range("I132").Value = 11930377
$ws.Range("J132").Value = 101520.73
$ws.Range("K132").Value = 35791131
$ws.Range("L132").Value = 304562.19
$ws.Range("M132").Value = -35788601
$ws.Range("N132").Value = -309622.19
# row 139
$ws.Range("H139").Value = 100000
$ws.Range("J139").Value = 100000
$ws.Range("L139").Value = 100000
$ws.Range("N139").Value = -110280

$ws = $wb.Worksheets.Item("BSM")
# row 141
$ws.Range("H141").Value = 75593.336
$ws.Range("J141").Value = 75593.336
$ws.Range("L141").Value = 75593.336
$ws.Range("N141").Value = -85953.336

$ws = $wb.Worksheets.Item("CUL")
# row 129
$ws.Range("H129").Value = 3474457.5
$ws.Range("I129").Value = 2152.5
$ws.Range("J129").Value = 4168918.5
$ws.Range("K129").Value = 6457.5
$ws.Range("L129").Value = 12506755.5
$ws.Range("M129").Value = -1457.5
$ws.Range("N129").Value = -12516755.5

$ws = $wb.Worksheets.Item("GSM")
# row 126
$ws.Range("H126").Value = 2407
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 3014
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 9042
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -13982

$ws = $wb.Worksheets.Item("LTW")
# row 16
$ws.Range("H16").Value = 2075.25
$ws.Range("I16").Value = 1900
$ws.Range("J16").Value = 2367.3333
$ws.Range("K16").Value = 1900
$ws.Range("L16").Value = 2367.3333
$ws.Range("M16").Value = -1730
$ws.Range("N16").Value = -2707.3333
# row 61
$ws.Range("H61").Value = 1787.4706
$ws.Range("I61").Value = 1812
$ws.Range("J61").Value = 1742.5
$ws.Range("K61").Value = 1812
$ws.Range("L61").Value = 1742.5
$ws.Range("M61").Value = -1610
$ws.Range("N61").Value = -2146.5
# row 93
$ws.Range("H93").Value = 1900
$ws.Range("I93").Value = 1900
$ws.Range("K93").Value = 1900
$ws.Range("M93").Value = -652
# row 100
$ws.Range("H100").Value = 1592.9286
$ws.Range("I100").Value = 1264.4286
$ws.Range("J100").Value = 1921.4286
$ws.Range("K100").Value = 1264.4286
$ws.Range("L100").Value = 1921.4286
$ws.Range("M100").Value = -723.4286
$ws.Range("N100").Value = -3003.4286
# row 113
$ws.Range("H113").Value = 1787.4706
$ws.Range("I113").Value = 1812
$ws.Range("J113").Value = 1742.5
$ws.Range("K113").Value = 1812
$ws.Range("L113").Value = 1742.5
$ws.Range("M113").Value = 358
$ws.Range("N113").Value = -6082.5
# row 132
$ws.Range("H132").Value = 25395.227
$ws.Range("I132").Value = 2645.2693
$ws.Range("K132").Value = 7935.8079
$ws.Range("M132").Value = -5405.8079

$ws = $wb.Worksheets.Item("WVR")
# row 62
$ws.Range("H62").Value = 3975.25
$ws.Range("I62").Value = 4002
$ws.Range("J62").Value = 3966.3333
$ws.Range("K62").Value = 4002
$ws.Range("L62").Value = 3966.3333
$ws.Range("M62").Value = -3378
$ws.Range("N62").Value = -5214.3333
# row 65
$ws.Range("H65").Value = 3975.25
$ws.Range("I65").Value = 4002
$ws.Range("J65").Value = 3966.3333
$ws.Range("K65").Value = 20010
$ws.Range("L65").Value = 19831.6665
$ws.Range("M65").Value = -16890
$ws.Range("N65").Value = -26071.6665
# row 81
$ws.Range("H81").Value = 2937.25
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 2937.25
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 5874.5
$ws.Range("N81").Value = -7996.5
$ws.Range("M81").ClearContents()
# row 84
$ws.Range("H84").Value = 2937.25
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 2937.25
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 29372.5
$ws.Range("N84").Value = -39980.5
$ws.Range("M84").ClearContents()
# row 132
$ws.Range("H132").Value = 48200.56
$ws.Range("I132").Value = 38172.074
$ws.Range("J132").Value = 65123.625
$ws.Range("K132").Value = 114516.222
$ws.Range("L132").Value = 195370.875
$ws.Range("M132").Value = -111986.222
$ws.Range("N132").Value = -200430.875
# row 136
$ws.Range("H136").Value = 31950.074
$ws.Range("I136").Value = 20726.176
$ws.Range("J136").Value = 67726.25
$ws.Range("K136").Value = 62178.528
$ws.Range("L136").Value = 203178.75
$ws.Range("M136").Value = -59628.528
$ws.Range("N136").Value = -208278.75
